$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value2 = '75.283.17'
$ws.Cells.Item(2, 5).Value2 = '  +7.37%  '

$ws.Cells.Item(3, 4).Value2 = '2.673.73'
$ws.Cells.Item(3, 5).Value2 = '  +9.07%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value2 = '187.45'
$ws.Cells.Item(5, 5).Value2 = '  +12.18%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value2 = '587.34'
$ws.Cells.Item(6, 5).Value2 = '  +3.35%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value2 = '0.998'
$ws.Cells.Item(7, 5).Value2 = '  -0.17%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value2 = '0.535'
$ws.Cells.Item(8, 5).Value2 = '  +4.01%  '

$ws.Cells.Item(9, 5).Value2 = '  +10.47%  '

$ws.Cells.Item(10, 4).Value2 = '2.673.02'
$ws.Cells.Item(10, 5).Value2 = '  +9.10%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value2 = '0.163'
$ws.Cells.Item(11, 5).Value2 = '  +1.37%  '

$ws.Cells.Item(12, 5).Value2 = '  +6.32%  '

$ws.Cells.Item(13, 5).Value2 = '  +0.28%  '

$ws.Cells.Item(14, 2).Value2 = 'WrappedBTC'
$ws.Cells.Item(14, 3).Value2 = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(14, 4).Value2 = '75.139.62'
$ws.Cells.Item(14, 5).Value2 = '  +7.21%  '

$ws.Cells.Item(15, 2).Value2 = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(15, 3).Value2 = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(15, 4).Value2 = '3.163.93'
$ws.Cells.Item(15, 5).Value2 = '  +8.95%  '

$ws.Cells.Item(16, 5).Value2 = '  +2.91%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value2 = '26.47'
$ws.Cells.Item(17, 5).Value2 = '  +9.30%  '

$ws.Cells.Item(18, 4).Value2 = '2.697.64'
$ws.Cells.Item(18, 5).Value2 = '  +9.98%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value2 = '9.16'
$ws.Cells.Item(19, 5).Value2 = '  +28.10%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value2 = '11.91'
$ws.Cells.Item(20, 5).Value2 = '  +9.39%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value2 = '371.78'
$ws.Cells.Item(21, 5).Value2 = '  +8.83%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value2 = '2.26'
$ws.Cells.Item(22, 5).Value2 = '  +12.37%  '

$ws.Cells.Item(23, 5).Value2 = '  +4.66%  '

$ws.Cells.Item(24, 5).Value2 = '  +3.55%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value2 = '0.999'
$ws.Cells.Item(25, 5).Value2 = '  -0.01%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value2 = '69.92'
$ws.Cells.Item(26, 5).Value2 = '  +5.29%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value2 = '4.13'
$ws.Cells.Item(27, 5).Value2 = '  +7.58%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value2 = '9.33'
$ws.Cells.Item(28, 5).Value2 = '  +9.46%  '

$ws.Cells.Item(29, 5).Value2 = '  +8.33%  '

$ws.Cells.Item(30, 5).Value2 = '  +0.46%  '

$ws.Cells.Item(31, 4).Value2 = '0.0₃0946'
$ws.Cells.Item(31, 5).Value2 = '  +9.83%  '

$ws.Cells.Item(32, 5).Value2 = '  +13.36%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value2 = '519.94'
$ws.Cells.Item(33, 5).Value2 = '  +12.41%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value2 = '7.68'
$ws.Cells.Item(34, 5).Value2 = '  +3.77%  '

$ws.Cells.Item(35, 5).Value2 = '  +7.67%  '

$ws.Cells.Item(36, 5).Value2 = '  -0.08%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value2 = '163.61'
$ws.Cells.Item(37, 5).Value2 = '  +2.69%  '

$ws.Cells.Item(38, 5).Value2 = '  +5.64%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value2 = '19.19'
$ws.Cells.Item(39, 5).Value2 = '  +5.18%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value2 = '19.36'
$ws.Cells.Item(40, 5).Value2 = '  +1.32%  '

$ws.Cells.Item(41, 5).Value2 = '  +0.02%  '

$ws.Cells.Item(42, 2).Value2 = 'Aave'
$ws.Cells.Item(42, 3).Value2 = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value2 = '170.29'
$ws.Cells.Item(42, 5).Value2 = '  +26.46%  '

$ws.Cells.Item(43, 2).Value2 = 'RenderToken'
$ws.Cells.Item(43, 3).Value2 = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value2 = '4.99'
$ws.Cells.Item(43, 5).Value2 = '  +12.62%  '

$ws.Cells.Item(44, 5).Value2 = '  +8.42%  '

$ws.Cells.Item(45, 5).Value2 = '  +9.63%  '

$ws.Cells.Item(46, 5).Value2 = '  +8.58%  '

$ws.Cells.Item(47, 5).Value2 = '  +11.52%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value2 = '39.03'
$ws.Cells.Item(48, 5).Value2 = '  +2.46%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value2 = '0.0848'
$ws.Cells.Item(49, 5).Value2 = '  +16.70%  '

$ws.Cells.Item(50, 5).Value2 = '  +7.02%  '

$ws.Cells.Item(51, 5).Value2 = '  +7.85%  '
